$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# "ajout ssr, had et psy" -- append a new record (row 46) to the format
# description table: type "c", position 182, nom "ZAD".
$ws.Range("B46").Value = "c"
$ws.Range("D46").Value = 182
$ws.Range("F46").Value = "ZAD"

# Leave the cursor just past the newly-added row, mirroring where the
# author's selection ended up after typing in the new data.
$ws.Range("D47").Select()
